# Update the test data sheet: the "TextToSearch" value used by the
# Google_Search_Test_2 test case changes from "Infosys" to "CGI".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Google_Search_Test_2")
$ws.Range("A3").Value = "CGI"
